$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two test rows (rows 3 and 4), keeping only the header and one data row
$ws.Rows("3:4").Delete()

# Update the remaining data row to describe a "NavigateTo" test case
$ws.Range("C2").Value = "NavigateTo"
$ws.Range("D2").Value = "NA"
$ws.Range("E2").ClearFormats()
$ws.Range("E2").Value = "https://www.google.com/"

# Turn the URL in E2 into a real hyperlink
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.google.com/")

# Resize columns C and E to fit their new content
$ws.Columns("C:C").ColumnWidth = 9.5
$ws.Columns("E:E").ColumnWidth = 22

# Update the active selection
$ws.Range("F10").Select() | Out-Null
